$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill column F ("OK") for rows that previously had no F value, matching column E pattern.
$okRows = @(3,4,5,6,7,8,9,11,12,13,14,15,16,17,18,19,20,21,22)
foreach ($r in $okRows) {
    $ws.Cells.Item($r, 6).Value = "OK"
}

# Adjust row heights
$ws.Rows.Item(4).RowHeight = 16.5
$ws.Rows.Item(8).RowHeight = 15.75

# New rows of test data (24 and 25)
$ws.Range("B24").Value = "Тормозит выполнение первой команды в пошаговом режиме"
$ws.Range("F24").Value = "NG"
$ws.Range("B25").Value = "При преобразовании кода условий в команды не учитывается  сравнивающая функция"
$ws.Range("F25").Value = "OK"

# Update the active selection to reflect the new view position
$ws.Activate()
$ws.Range("F26").Select()
